# The RTM table lists requirement rows LTT_01..LTT_18 in columns A:E.
# Row 18 (LTT_17 / 7.1.17) was mapped to "Book_ticket" / "Test_case_2.1",
# while row 19 (LTT_18 / 7.1.18) was mapped to "Make journey" / "Test_case_2.2".
# The edit removes the separate "LTT_18" requirement row, and re-points the
# "LTT_17" row's Code/UT mapping (columns C:D) to what used to be row 19's
# values ("Make journey" / "Test_case_2.2"), dropping the now-unused
# "Book_ticket" / "Test_case_2.1" / "7.1.18" / "LTT_18" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull row 19's Code Mapping / UT Mapping values up into row 18, overwriting
# the old "Book_ticket" / "Test_case_2.1" entries.
$ws.Range("C18").Value2 = $ws.Range("C19").Value2
$ws.Range("D18").Value2 = $ws.Range("D19").Value2

# Remove the now-redundant last row (LTT_18 / 7.1.18) entirely, shifting the
# used range back up to A1:E18.
$ws.Rows("19").Delete()

# Match the saved cursor position left behind by this edit.
$ws.Range("D19").Select()
